$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 1.39
$ws.Range("F3").Value = 1.19
$ws.Range("G3").Value = 0.64

# Row 4
$ws.Range("D4").Value = 1.32

# Row 5
$ws.Range("C5").Value = 1.38
$ws.Range("G5").Value = 0.77

# Row 7
$ws.Range("C7").Value = 2.2
$ws.Range("E7").Value = 1.87
